$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2010 and 2010-18")

# --- Row 75: K75 and Q75 lose their yellow highlight fill -----------------
# (ClearFormats + re-apply the 0.00 number format reproduces the plain,
#  un-highlighted style used elsewhere on this sheet)
$ws.Range("K75").ClearFormats()
$ws.Range("K75").Value = 676.83127177777783
$ws.Range("K75").NumberFormat = "0.00"

$ws.Range("Q75").ClearFormats()
$ws.Range("Q75").Value = 0.14695311111111112
$ws.Range("Q75").NumberFormat = "0.00"

# --- Row 76: new scenario row "Baseline_1979-current C286" ----------------
$ws.Range("A76").Value = "CW3M"
$ws.Range("B76").Value = "Baseline_1979-current C286"
$ws.Range("C76").Value = "2010-18"

$ws.Range("D76").Value = 1112.7018771111111
$ws.Range("D76").NumberFormat = "0.00"

$ws.Range("E76").Value = 1763.5263265555557
$ws.Range("E76").NumberFormat = "0.00"
$ws.Range("E76").Interior.Color = 65535

$ws.Range("F76").Value = 1.1070731111111112
$ws.Range("F76").NumberFormat = "0.00"

$ws.Range("G76").Value = 295.25833466666666
$ws.Range("G76").NumberFormat = "0.00"

$ws.Range("H76").Value = 9.775355222222224
$ws.Range("H76").NumberFormat = "0.00"

$ws.Range("I76").Value = 6.5172971111111115
$ws.Range("I76").NumberFormat = "0.00"

$ws.Range("J76").Value = 8.145128999999999
$ws.Range("J76").NumberFormat = "0.00"

$ws.Range("K76").Value = 686.8172538888889
$ws.Range("K76").NumberFormat = "0.00"
$ws.Range("K76").Interior.Color = 65535

$ws.Range("L76").Value = 60.018756111111117
$ws.Range("L76").NumberFormat = "0.00"
$ws.Range("L76").Interior.Color = 65535

$ws.Range("M76").Value = 1361.733412
$ws.Range("M76").NumberFormat = "0.00"
$ws.Range("M76").Interior.Color = 65535

$ws.Range("N76").Value = 1072.4035372222222
$ws.Range("N76").NumberFormat = "0.00"
$ws.Range("N76").Interior.Color = 65535

$ws.Range("O76").Value = 6355.1079644444444
$ws.Range("O76").NumberFormat = "0"
$ws.Range("O76").Interior.Color = 65535

$ws.Range("P76").Value = 27227.338324888889
$ws.Range("P76").NumberFormat = "0"

$ws.Range("Q76").ClearFormats()
$ws.Range("Q76").Value = 0.23182455555555562
$ws.Range("Q76").NumberFormat = "0.00"

$ws.Range("R76").Value = 0.000041777777777777767
$ws.Range("R76").NumberFormat = "0.000000"

# Reproduce the new active selection noted in the saved file.
$ws.Range("K75").Select()
